# Updates cryptocurrency price/volume data per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '35.274.83'
$ws.Range("E2").Value = '  +0.56%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.896.01'
$ws.Range("E3").Value = '  +2.49%  '

# Row 4
$ws.Range("E4").Value = '  +0.09%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.96'
$ws.Range("E5").Value = '  +2.33%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.650'
$ws.Range("E6").Value = '  +4.76%  '

# Row 7
$ws.Range("E7").Value = '  +0.25%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '41.44'
$ws.Range("E8").Value = '  -2.47%  '

# Row 9
$ws.Range("E9").Value = '  +4.63%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '50.30'
$ws.Range("E10").Value = '  +8.25%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0707'
$ws.Range("E11").Value = '  +2.39%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0995'
$ws.Range("E12").Value = '  +0.67%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.174.56'
$ws.Range("E13").Value = '  +2.67%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '12.06'
$ws.Range("E14").Value = '  +6.02%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.691'
$ws.Range("E15").Value = '  +2.47%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.902.15'
$ws.Range("E16").Value = '  +2.85%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.81'
$ws.Range("E17").Value = '  +1.25%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '35.332.77'
$ws.Range("E18").Value = '  +0.84%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '71.27'
$ws.Range("E19").Value = '  +1.89%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0813'
$ws.Range("E20").Value = '  +2.80%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '241.14'
$ws.Range("E21").Value = '  +0.30%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.47'
$ws.Range("E22").Value = '  +2.72%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.70'
$ws.Range("E23").Value = '  -0.90%  '

# Row 24
$ws.Range("E24").Value = '  +0.13%  '

# Row 25
$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.38'
$ws.Range("E25").Value = '  +30.11%  '

# Row 26
$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.30'
$ws.Range("E26").Value = '  +1.36%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '170.02'
$ws.Range("E27").Value = '  +0.16%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.31'
$ws.Range("E28").Value = '  +3.77%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.14'
$ws.Range("E29").Value = '  +3.13%  '

# Row 30
$ws.Range("E30").Value = '  +1.90%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.10'
$ws.Range("E31").Value = '  +2.88%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0559'
$ws.Range("E32").Value = '  +1.36%  '

# Row 33
$ws.Range("E33").Value = '  +0.00%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.922'
$ws.Range("E34").Value = '  +15.89%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.09'
$ws.Range("E35").Value = '  +2.05%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.74'
$ws.Range("E36").Value = '  +1.42%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.02'
$ws.Range("E37").Value = '  +1.34%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.32'
$ws.Range("E38").Value = '  +1.78%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0209'
$ws.Range("E39").Value = '  +3.85%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.08'
$ws.Range("E40").Value = '  +1.34%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0631'
$ws.Range("E41").Value = '  +14.23%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '15.83'
$ws.Range("E42").Value = '  +6.56%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '88.85'
$ws.Range("E43").Value = '  -1.40%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.334.40'
$ws.Range("E44").Value = '  -0.53%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.34'
$ws.Range("E45").Value = '  +1.83%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '46.81'
$ws.Range("E46").Value = '  +36.19%  '

# Row 47
$ws.Range("E47").Value = '  -1.72%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.77'
$ws.Range("E48").Value = '  +1.48%  '

# Row 49
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.48'
$ws.Range("E49").Value = '  +0.34%  '

# Row 50
$ws.Range("B50").Value = 'Gas'
$ws.Range("C50").Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '11.21'
$ws.Range("E50").Value = '  -13.64%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.083.91'
$ws.Range("E51").Value = '  +2.42%  '
